$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4130
$ws.Range("I64").Value = 4299
$ws.Range("J64").Value = 3976.3635
$ws.Range("K64").Value = 4299
$ws.Range("L64").Value = 3976.3635
$ws.Range("M64").Value = -4051
$ws.Range("N64").Value = -4472.363499999999

$ws.Range("H67").Value = 4130
$ws.Range("I67").Value = 4299
$ws.Range("J67").Value = 3976.3635
$ws.Range("K67").Value = 4299
$ws.Range("L67").Value = 3976.3635
$ws.Range("M67").Value = -3441
$ws.Range("N67").Value = -5692.363499999999

$ws.Range("H116").Value = 3412.7273
$ws.Range("I116").Value = 3728.5715
$ws.Range("J116").Value = 2860
$ws.Range("K116").Value = 3728.5715
$ws.Range("L116").Value = 2860
$ws.Range("M116").Value = -286.5715
$ws.Range("N116").Value = -9744

$ws.Range("H132").Value = 8777154
$ws.Range("I132").Value = 11911013
$ws.Range("K132").Value = 35733039
$ws.Range("M132").Value = -35730509

$ws.Range("H135").Value = 32258752
$ws.Range("I135").Value = 330.91666
$ws.Range("J135").Value = 142859060
$ws.Range("K135").Value = 2978.24994
$ws.Range("L135").Value = 1285731540
$ws.Range("M135").Value = -443.2499399999997
$ws.Range("N135").Value = -1285736610

$ws.Range("H138").Value = 525552.5600000001
$ws.Range("I138").Value = 1430.037
$ws.Range("J138").Value = 778254.5600000001
$ws.Range("K138").Value = 4290.111
$ws.Range("L138").Value = 2334763.68
$ws.Range("M138").Value = 849.8890000000001
$ws.Range("N138").Value = -2345043.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2934.674
$ws.Range("I32").Value = 2622.1013
$ws.Range("K32").Value = 2622.1013
$ws.Range("M32").Value = -2335.1013

$ws.Range("H45").Value = 1099.5518
$ws.Range("I45").Value = 1075.1177
$ws.Range("J45").Value = 1134.1666
$ws.Range("K45").Value = 1075.1177
$ws.Range("L45").Value = 1134.1666
$ws.Range("M45").Value = -698.1177
$ws.Range("N45").Value = -1888.1666

$ws.Range("H61").Value = 1475.7241
$ws.Range("I61").Value = 1304.8334
$ws.Range("K61").Value = 1304.8334
$ws.Range("M61").Value = -1092.8334

$ws.Range("H108").Value = 9683
$ws.Range("J108").Value = 9683
$ws.Range("L108").Value = 9683
$ws.Range("N108").Value = -17363

$ws.Range("H122").Value = 810.7895
$ws.Range("I122").Value = 800.2778
$ws.Range("K122").Value = 2400.8334
$ws.Range("M122").Value = 49.16660000000002

$ws.Range("H132").Value = 2568.6428
$ws.Range("I132").Value = 2163.5833
$ws.Range("K132").Value = 6490.749899999999
$ws.Range("M132").Value = -3960.749899999999

$ws.Range("H136").Value = 1475.7241
$ws.Range("I136").Value = 1304.8334
$ws.Range("K136").Value = 3914.5002
$ws.Range("M136").Value = -1364.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1374.9131
$ws.Range("I20").Value = 918.55554
$ws.Range("K20").Value = 918.55554
$ws.Range("M20").Value = -671.55554

$ws.Range("H105").Value = 250003300
$ws.Range("I105").Value = 250003300
$ws.Range("K105").Value = 250003300
$ws.Range("M105").Value = -250001553

$ws.Range("H132").Value = 45374.625
$ws.Range("J132").Value = 45374.625
$ws.Range("L132").Value = 45374.625
$ws.Range("N132").Value = -55494.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1271.1632
$ws.Range("I31").Value = 1240.1489
$ws.Range("K31").Value = 1240.1489
$ws.Range("M31").Value = -945.1488999999999

$ws.Range("H34").Value = 1271.1632
$ws.Range("I34").Value = 1240.1489
$ws.Range("K34").Value = 1240.1489
$ws.Range("M34").Value = -1038.1489

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 346.75
$ws.Range("I92").Value = 624.5
$ws.Range("J92").Value = 254.16667
$ws.Range("K92").Value = 1873.5
$ws.Range("L92").Value = 762.50001
$ws.Range("M92").Value = -625.5
$ws.Range("N92").Value = -3258.50001

$ws.Range("H117").Value = 555.3333
$ws.Range("I117").Value = 566.6667
$ws.Range("J117").Value = 544
$ws.Range("K117").Value = 1700.0001
$ws.Range("L117").Value = 1632
$ws.Range("M117").Value = 1741.9999
$ws.Range("N117").Value = -8516

$ws.Range("H131").Value = 18870864
$ws.Range("I131").Value = 333333660
$ws.Range("J131").Value = 3095.08
$ws.Range("K131").Value = 1000000980
$ws.Range("L131").Value = 9285.24
$ws.Range("M131").Value = -999995940
$ws.Range("N131").Value = -19365.24

$ws.Range("H132").Value = 1639.2778
$ws.Range("I132").Value = 1050.5834
$ws.Range("J132").Value = 2816.6667
$ws.Range("K132").Value = 9455.250599999999
$ws.Range("L132").Value = 25350.0003
$ws.Range("M132").Value = -6925.250599999999
$ws.Range("N132").Value = -30410.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2714.3572
$ws.Range("I132").Value = 2453.8462
$ws.Range("J132").Value = 2940.1333
$ws.Range("K132").Value = 7361.5386
$ws.Range("L132").Value = 8820.3999
$ws.Range("M132").Value = -4831.5386
$ws.Range("N132").Value = -13880.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2039
$ws.Range("I7").Value = 1998.8889
$ws.Range("J7").Value = 2400
$ws.Range("K7").Value = 1998.8889
$ws.Range("L7").Value = 2400
$ws.Range("M7").Value = -1886.8889
$ws.Range("N7").Value = -2624

$ws.Range("H40").Value = 3335.9092
$ws.Range("I40").Value = 2232.0833
$ws.Range("J40").Value = 4660.5
$ws.Range("K40").Value = 2232.0833
$ws.Range("L40").Value = 4660.5
$ws.Range("M40").Value = -2096.0833
$ws.Range("N40").Value = -4932.5

$ws.Range("H46").Value = 5929.9
$ws.Range("I46").Value = 2999.5
$ws.Range("J46").Value = 6662.5
$ws.Range("K46").Value = 2999.5
$ws.Range("L46").Value = 6662.5
$ws.Range("M46").Value = -2811.5
$ws.Range("N46").Value = -7038.5

$ws.Range("H82").Value = 2314.7693
$ws.Range("I82").Value = 2500
$ws.Range("K82").Value = 2500
$ws.Range("M82").Value = -2139

$ws.Range("H85").Value = 2314.7693
$ws.Range("I85").Value = 2500
$ws.Range("K85").Value = 2500
$ws.Range("M85").Value = -1252

$ws.Range("H122").Value = 19241146
$ws.Range("I122").Value = 41685516
$ws.Range("J122").Value = 3114.2856
$ws.Range("K122").Value = 125056548
$ws.Range("L122").Value = 9342.856800000001
$ws.Range("M122").Value = -125054098
$ws.Range("N122").Value = -14242.8568

$ws.Range("H126").Value = 2039
$ws.Range("I126").Value = 1998.8889
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5996.6667
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -3526.6667
$ws.Range("N126").Value = -12140

$ws.Range("H132").Value = 28355
$ws.Range("I132").Value = 1583.5264
$ws.Range("J132").Value = 55126.473
$ws.Range("K132").Value = 4750.5792
$ws.Range("L132").Value = 165379.419
$ws.Range("M132").Value = -2220.5792
$ws.Range("N132").Value = -170439.419

$ws.Range("H136").Value = 5495.8696
$ws.Range("I136").Value = 6605
$ws.Range("J136").Value = 1503
$ws.Range("K136").Value = 19815
$ws.Range("L136").Value = 4509
$ws.Range("M136").Value = -17265
$ws.Range("N136").Value = -9609

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10377
$ws.Range("J41").Value = 10377
$ws.Range("L41").Value = 10377
$ws.Range("N41").Value = -11157

$ws.Range("H81").Value = 399.33334
$ws.Range("I81").Value = 399.33334
$ws.Range("K81").Value = 798.66668
$ws.Range("M81").Value = 262.33332

$ws.Range("H84").Value = 399.33334
$ws.Range("I84").Value = 399.33334
$ws.Range("K84").Value = 3993.3334
$ws.Range("M84").Value = 1310.6666

$ws.Range("H132").Value = 3246.3142
$ws.Range("I132").Value = 3345.36
$ws.Range("J132").Value = 2998.7
$ws.Range("K132").Value = 10036.08
$ws.Range("L132").Value = 8996.099999999999
$ws.Range("M132").Value = -7506.08
$ws.Range("N132").Value = -14056.1

$ws.Range("H136").Value = 606.2105
$ws.Range("I136").Value = 457.29413
$ws.Range("J136").Value = 1872
$ws.Range("K136").Value = 1371.88239
$ws.Range("L136").Value = 5616
$ws.Range("M136").Value = 1178.11761
$ws.Range("N136").Value = -10716
